# Applies the updated attribution values to row 2 of the active sheet
# (relative-direction update per commit message), cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = [double]"-0"  # A2
$ws.Cells.Item(2, 2).Value = [double]"-0.08608604857666213"  # B2
$ws.Cells.Item(2, 4).Value = [double]"0.2707430686247467"  # D2
$ws.Cells.Item(2, 5).Value = [double]"0.004001273688284372"  # E2
$ws.Cells.Item(2, 7).Value = [double]"0"  # G2
$ws.Cells.Item(2, 9).Value = [double]"-0"  # I2
$ws.Cells.Item(2, 10).Value = [double]"-0"  # J2
$ws.Cells.Item(2, 11).Value = [double]"-0.05711919423029652"  # K2
$ws.Cells.Item(2, 12).Value = [double]"-0"  # L2
$ws.Cells.Item(2, 13).Value = [double]"0.2424185176332642"  # M2
$ws.Cells.Item(2, 14).Value = [double]"0.01439279748655964"  # N2
$ws.Cells.Item(2, 15).Value = [double]"0"  # O2
$ws.Cells.Item(2, 19).Value = [double]"0"  # S2
$ws.Cells.Item(2, 20).Value = [double]"-0.1039938171523975"  # T2
$ws.Cells.Item(2, 22).Value = [double]"0.01705956103150915"  # V2
$ws.Cells.Item(2, 23).Value = [double]"-0.01078891192103048"  # W2
$ws.Cells.Item(2, 25).Value = [double]"-0"  # Y2
$ws.Cells.Item(2, 26).Value = [double]"-0"  # Z2
$ws.Cells.Item(2, 28).Value = [double]"0"  # AB2
$ws.Cells.Item(2, 29).Value = [double]"-3.470842752942556e-11"  # AC2
$ws.Cells.Item(2, 30).Value = [double]"0"  # AD2
$ws.Cells.Item(2, 31).Value = [double]"0.006287610896074329"  # AE2
$ws.Cells.Item(2, 32).Value = [double]"-1.992104022987901e-12"  # AF2
$ws.Cells.Item(2, 33).Value = [double]"-0"  # AG2
$ws.Cells.Item(2, 35).Value = [double]"-0"  # AI2
$ws.Cells.Item(2, 36).Value = [double]"0"  # AJ2
$ws.Cells.Item(2, 37).Value = [double]"0"  # AK2
$ws.Cells.Item(2, 38).Value = [double]"-0.04586640419228727"  # AL2
$ws.Cells.Item(2, 39).Value = [double]"0"  # AM2
$ws.Cells.Item(2, 40).Value = [double]"0.04249581251515201"  # AN2
$ws.Cells.Item(2, 41).Value = [double]"0.07080012785383533"  # AO2
$ws.Cells.Item(2, 43).Value = [double]"0"  # AQ2
$ws.Cells.Item(2, 44).Value = [double]"-0"  # AR2
$ws.Cells.Item(2, 45).Value = [double]"-0"  # AS2
$ws.Cells.Item(2, 46).Value = [double]"0"  # AT2
$ws.Cells.Item(2, 47).Value = [double]"-0.1675629633249888"  # AU2
$ws.Cells.Item(2, 49).Value = [double]"0.09601824582853803"  # AW2
$ws.Cells.Item(2, 50).Value = [double]"-0.01667735686670889"  # AX2
$ws.Cells.Item(2, 51).Value = [double]"-0"  # AY2
$ws.Cells.Item(2, 52).Value = [double]"0"  # AZ2
$ws.Cells.Item(2, 54).Value = [double]"-0"  # BB2
$ws.Cells.Item(2, 55).Value = [double]"-0"  # BC2
$ws.Cells.Item(2, 56).Value = [double]"-0.007563276064319164"  # BD2
$ws.Cells.Item(2, 58).Value = [double]"0.1049887927948767"  # BF2
$ws.Cells.Item(2, 59).Value = [double]"0.004555816567386153"  # BG2
$ws.Cells.Item(2, 62).Value = [double]"-0"  # BJ2
$ws.Cells.Item(2, 64).Value = [double]"0"  # BL2
$ws.Cells.Item(2, 65).Value = [double]"0.01855919159624354"  # BM2
$ws.Cells.Item(2, 67).Value = [double]"-0.03870817360056667"  # BO2
$ws.Cells.Item(2, 68).Value = [double]"-0.08014909847657475"  # BP2
$ws.Cells.Item(2, 73).Value = [double]"0"  # BU2
$ws.Cells.Item(2, 74).Value = [double]"-0.06240362818182101"  # BV2
$ws.Cells.Item(2, 75).Value = [double]"0"  # BW2
$ws.Cells.Item(2, 76).Value = [double]"0.02238961253711967"  # BX2
$ws.Cells.Item(2, 77).Value = [double]"-0.02819281016782701"  # BY2
$ws.Cells.Item(2, 78).Value = [double]"-0"  # BZ2
$ws.Cells.Item(2, 80).Value = [double]"0"  # CB2
$ws.Cells.Item(2, 82).Value = [double]"-0"  # CD2
$ws.Cells.Item(2, 83).Value = [double]"0.03017968464745394"  # CE2
$ws.Cells.Item(2, 85).Value = [double]"-0.04373135439866592"  # CG2
$ws.Cells.Item(2, 86).Value = [double]"0.01555492970330826"  # CH2
$ws.Cells.Item(2, 88).Value = [double]"-0"  # CJ2
$ws.Cells.Item(2, 91).Value = [double]"-0"  # CM2
$ws.Cells.Item(2, 92).Value = [double]"-0.01683645637809813"  # CN2
$ws.Cells.Item(2, 93).Value = [double]"-0"  # CO2
$ws.Cells.Item(2, 94).Value = [double]"0.02927517608387242"  # CP2
$ws.Cells.Item(2, 95).Value = [double]"0.05732951286099017"  # CQ2
$ws.Cells.Item(2, 96).Value = [double]"0"  # CR2
$ws.Cells.Item(2, 98).Value = [double]"0"  # CT2
$ws.Cells.Item(2, 99).Value = [double]"-0"  # CU2
$ws.Cells.Item(2, 100).Value = [double]"-0"  # CV2
$ws.Cells.Item(2, 101).Value = [double]"0.04485442205455682"  # CW2
$ws.Cells.Item(2, 103).Value = [double]"-0.04610213201090916"  # CY2
$ws.Cells.Item(2, 104).Value = [double]"0.007004487323011817"  # CZ2
$ws.Cells.Item(2, 105).Value = [double]"0"  # DA2
$ws.Cells.Item(2, 109).Value = [double]"-0"  # DE2
$ws.Cells.Item(2, 110).Value = [double]"0.03663010099907456"  # DF2
$ws.Cells.Item(2, 112).Value = [double]"0.01496390616943506"  # DH2
$ws.Cells.Item(2, 113).Value = [double]"0.03101452545121341"  # DI2
$ws.Cells.Item(2, 114).Value = [double]"0"  # DJ2
$ws.Cells.Item(2, 115).Value = [double]"-0"  # DK2
$ws.Cells.Item(2, 116).Value = [double]"-0"  # DL2
$ws.Cells.Item(2, 118).Value = [double]"0"  # DN2
$ws.Cells.Item(2, 119).Value = [double]"-0.02091649284860365"  # DO2
$ws.Cells.Item(2, 121).Value = [double]"0.04285491153331428"  # DQ2
$ws.Cells.Item(2, 122).Value = [double]"-0.03731491465691677"  # DR2
$ws.Cells.Item(2, 123).Value = [double]"-0"  # DS2
$ws.Cells.Item(2, 125).Value = [double]"-0"  # DU2
$ws.Cells.Item(2, 127).Value = [double]"0"  # DW2
$ws.Cells.Item(2, 128).Value = [double]"-0.05320440195189079"  # DX2
$ws.Cells.Item(2, 129).Value = [double]"-0"  # DY2
$ws.Cells.Item(2, 130).Value = [double]"-0.002892739116739693"  # DZ2
$ws.Cells.Item(2, 131).Value = [double]"-0.03046001104498006"  # EA2
$ws.Cells.Item(2, 132).Value = [double]"0"  # EB2
$ws.Cells.Item(2, 136).Value = [double]"-0"  # EF2
$ws.Cells.Item(2, 137).Value = [double]"0.04214883444383082"  # EG2
$ws.Cells.Item(2, 138).Value = [double]"0"  # EH2
$ws.Cells.Item(2, 139).Value = [double]"0.09694243965812084"  # EI2
$ws.Cells.Item(2, 140).Value = [double]"-0.02535434639816732"  # EJ2
$ws.Cells.Item(2, 141).Value = [double]"0"  # EK2
$ws.Cells.Item(2, 145).Value = [double]"0"  # EO2
$ws.Cells.Item(2, 146).Value = [double]"0.0523474947198154"  # EP2
$ws.Cells.Item(2, 147).Value = [double]"0"  # EQ2
$ws.Cells.Item(2, 148).Value = [double]"-0.0505857322256974"  # ER2
$ws.Cells.Item(2, 149).Value = [double]"0.02099614862398629"  # ES2
$ws.Cells.Item(2, 150).Value = [double]"0"  # ET2
$ws.Cells.Item(2, 152).Value = [double]"0"  # EV2
$ws.Cells.Item(2, 154).Value = [double]"0"  # EX2
$ws.Cells.Item(2, 155).Value = [double]"0.03977438091833668"  # EY2
$ws.Cells.Item(2, 156).Value = [double]"0"  # EZ2
$ws.Cells.Item(2, 157).Value = [double]"-0.03758484801388312"  # FA2
$ws.Cells.Item(2, 158).Value = [double]"0.01722156883334725"  # FB2
$ws.Cells.Item(2, 160).Value = [double]"-0"  # FD2
$ws.Cells.Item(2, 162).Value = [double]"-0"  # FF2
$ws.Cells.Item(2, 163).Value = [double]"-0"  # FG2
$ws.Cells.Item(2, 164).Value = [double]"-0.006986232390577869"  # FH2
$ws.Cells.Item(2, 166).Value = [double]"-0.0160849485541347"  # FJ2
$ws.Cells.Item(2, 167).Value = [double]"0.02278434274117715"  # FK2
$ws.Cells.Item(2, 168).Value = [double]"-0"  # FL2
$ws.Cells.Item(2, 170).Value = [double]"-0"  # FN2
$ws.Cells.Item(2, 172).Value = [double]"-0"  # FP2
$ws.Cells.Item(2, 173).Value = [double]"-0.009903314644151768"  # FQ2
$ws.Cells.Item(2, 174).Value = [double]"-0"  # FR2
$ws.Cells.Item(2, 175).Value = [double]"0.01161490472161882"  # FS2
$ws.Cells.Item(2, 176).Value = [double]"-0.00335253712351679"  # FT2
$ws.Cells.Item(2, 177).Value = [double]"0"  # FU2
$ws.Cells.Item(2, 178).Value = [double]"-0"  # FV2
$ws.Cells.Item(2, 179).Value = [double]"-0"  # FW2
$ws.Cells.Item(2, 181).Value = [double]"0"  # FY2
$ws.Cells.Item(2, 182).Value = [double]"-0.03276389726031452"  # FZ2
$ws.Cells.Item(2, 184).Value = [double]"0.02249852629173317"  # GB2
$ws.Cells.Item(2, 186).Value = [double]"0"  # GD2
$ws.Cells.Item(2, 187).Value = [double]"-0"  # GE2
